# Update countries & provincias Spain
# Applies the covid-data refresh described by the commit:
#  - Costa Rica now ranks ahead of China (row 50/51 swap positions)
#  - Guayana Francesa now ranks ahead of Tayikistan (row 107/108 swap positions)
#  - Refreshed totals for several country rows
#  - Updated "last updated" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Cells.Item(4, 2).Value2 = 7875633
$ws.Cells.Item(4, 3).Value2 = 41870
$ws.Cells.Item(4, 4).Value2 = 5050098
$ws.Cells.Item(4, 5).Value2 = 2607218
$ws.Cells.Item(4, 7).Value2 = 579
$ws.Cells.Item(4, 8).Value2 = 218317

# --- Row 5: India ---
$ws.Cells.Item(5, 2).Value2 = 6977008
$ws.Cells.Item(5, 3).Value2 = 73196
$ws.Cells.Item(5, 4).Value2 = 5985505
$ws.Cells.Item(5, 5).Value2 = 884053
$ws.Cells.Item(5, 7).Value2 = 929
$ws.Cells.Item(5, 8).Value2 = 107450

# --- Row 25: Alemania ---
$ws.Cells.Item(25, 2).Value2 = 320214
$ws.Cells.Item(25, 3).Value2 = 4700
$ws.Cells.Item(25, 5).Value2 = 41032
$ws.Cells.Item(25, 7).Value2 = 15
$ws.Cells.Item(25, 8).Value2 = 9682

# --- Row 50: now Costa Rica (was China) ---
$ws.Cells.Item(50, 1).Value2 = "Costa Rica"
$ws.Cells.Item(50, 2).Value2 = 86053
$ws.Cells.Item(50, 3).Value2 = 1225
$ws.Cells.Item(50, 4).Value2 = 52327
$ws.Cells.Item(50, 5).Value2 = 32671
$ws.Cells.Item(50, 7).Value2 = 15
$ws.Cells.Item(50, 8).Value2 = 1055

# --- Row 51: now China (was Costa Rica) ---
$ws.Cells.Item(51, 1).Value2 = "China"
$ws.Cells.Item(51, 2).Value2 = 85521
$ws.Cells.Item(51, 3).Value2 = 21
$ws.Cells.Item(51, 4).Value2 = 80681
$ws.Cells.Item(51, 5).Value2 = 206
$ws.Cells.Item(51, 7).Value2 = 0
$ws.Cells.Item(51, 8).Value2 = 4634

# --- Row 107: now Guayana Francesa (was Tayikistan) ---
$ws.Cells.Item(107, 1).Value2 = "Guayana Francesa"
$ws.Cells.Item(107, 2).Value2 = 10144
$ws.Cells.Item(107, 3).Value2 = 16
$ws.Cells.Item(107, 4).Value2 = 9810
$ws.Cells.Item(107, 5).Value2 = 265
$ws.Cells.Item(107, 7).Value2 = 0
$ws.Cells.Item(107, 8).Value2 = 69

# --- Row 108: now Tayikistan (was Guayana Francesa) ---
$ws.Cells.Item(108, 1).Value2 = "Tayikistan"
$ws.Cells.Item(108, 2).Value2 = 10137
$ws.Cells.Item(108, 3).Value2 = 40
$ws.Cells.Item(108, 4).Value2 = 8959
$ws.Cells.Item(108, 5).Value2 = 1099
$ws.Cells.Item(108, 7).Value2 = 1
$ws.Cells.Item(108, 8).Value2 = 79

# --- Row 120: Angola ---
$ws.Cells.Item(120, 2).Value2 = 6031
$ws.Cells.Item(120, 3).Value2 = 73
$ws.Cells.Item(120, 4).Value2 = 2685
$ws.Cells.Item(120, 5).Value2 = 3134
$ws.Cells.Item(120, 7).Value2 = 4
$ws.Cells.Item(120, 8).Value2 = 212

# --- Row 166: Republica del Chad ---
$ws.Cells.Item(166, 2).Value2 = 1274
$ws.Cells.Item(166, 3).Value2 = 12
$ws.Cells.Item(166, 4).Value2 = 1102
$ws.Cells.Item(166, 5).Value2 = 82
$ws.Cells.Item(166, 7).Value2 = 1
$ws.Cells.Item(166, 8).Value2 = 90

# --- Row 191: Barbados ---
$ws.Cells.Item(191, 2).Value2 = 204
$ws.Cells.Item(191, 3).Value2 = 1
$ws.Cells.Item(191, 5).Value2 = 15

# --- Row 198: Islas Virgenes Britanicas ---
$ws.Cells.Item(198, 4).Value2 = 70
$ws.Cells.Item(198, 5).Value2 = 0

# --- Updated timestamp banner ---
$ws.Cells.Item(1, 1).Value2 = "Datos actualizados a 9 de Octubre de 2020 a las 22:07"
